# Team Log workbook update — "Filling up TeamLog - Kumai"
# Replaces the old "navbar" entry (row 16) with Kumai's five new Sprint 2
# log entries, and updates row heights that Excel recomputed as a result
# of the content reflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Overwrite the old row 16 ("Worked on turning the navbar into a
#    reuseable component..." / Gabriel and Kumai / 7h) with the first
#    of Kumai's new entries. The cell style (s="5") is preserved as-is
#    since we are only touching the value, not the formatting.
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Creating Sql querries and associated JavaScript functions for frontend-backend communications"
$ws.Range("B16").Value = "I had communicated with Gabriel and Yousef to make sure the needed data, filer parameters, and all database related functions are available for the pages development."
$ws.Range("C16").Value = "Steven and Kumai"
$ws.Range("D16").Value = "8h"

# ---------------------------------------------------------------------
# 2. Write four more new rows below (17:20) — nothing existed past row
#    16 before, so this simply extends the sheet; no row-shifting is
#    required. Fill in the values first, then copy row 16's formatting
#    (style s="5") down onto them.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Working on the property search front-end back-end integration and improving UI"
$ws.Range("B17").Value = "There are still work to be done to improve the user interface. For now, the user can easily search for properties by the area, maximum price, and if it is for sale or for rent."
$ws.Range("C17").Value = "Gabriel and Kumai"
$ws.Range("D17").Value = "5h"

$ws.Range("A18").Value = "Initializing CI/CD Pipeline on GitHub and making a basic testing architecture for the project"
$ws.Range("B18").Value = "It is difficult to perform test cases when the whole project is locally hosted. The test cases are done locally and the work perfectly. However, they do not work on GitHub since there is no way to connect to the local database. Hence, I have setup hard coded data for now, but this can be easily changed once a better database enviornment is set. "
$ws.Range("C18").Value = "Kumai"
$ws.Range("D18").Value = "8h"

$ws.Range("A19").Value = "Setting up test cases for the project"
$ws.Range("B19").Value = "There is only one php file to be tested. Other files are all JavaScript files. However, there are too many sql querries related functions that need to be tested."
$ws.Range("C19").Value = "Kumai"
$ws.Range("D19").Value = "5h"

$ws.Range("A20").Value = "Peer reviewing and support"
$ws.Range("B20").Value = "Helping other teams with the setup and debugging. "
$ws.Range("C20").Value = "Kumai"
$ws.Range("D20").Value = "6h"

$ws.Range("A16:D16").Copy()
$ws.Range("A17:D20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Row heights reflowed across the whole sheet once the new content
#    was added (rows 4 & 5 keep their existing custom height).
# ---------------------------------------------------------------------
$ws.Rows(1).RowHeight = 18.75
$ws.Rows(2).RowHeight = 18.75
$ws.Rows(3).RowHeight = 90
$ws.Rows(6).RowHeight = 60
$ws.Rows(7).RowHeight = 30
$ws.Rows(8).RowHeight = 30
$ws.Rows(9).RowHeight = 30
$ws.Rows(10).RowHeight = 90
$ws.Rows(11).RowHeight = 45
$ws.Rows(12).RowHeight = 30
$ws.Rows(13).RowHeight = 18.75
$ws.Rows(14).RowHeight = 30
$ws.Rows(15).RowHeight = 75
$ws.Rows(16).RowHeight = 75
$ws.Rows(17).RowHeight = 75
$ws.Rows(18).RowHeight = 135
$ws.Rows(19).RowHeight = 60
$ws.Rows(20).RowHeight = 30

# ---------------------------------------------------------------------
# 4. Leave the selection where the editor (Kumai) left off.
# ---------------------------------------------------------------------
$ws.Range("C24").Select()
